# Update "Inscricoes" sheet numeric cells per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @{
    "E2" = 60;
    "F2" = 42;
    "H2" = 54;
    "E7" = 14;
    "F7" = 5;
    "H7" = 6;
    "E8" = 18;
    "E9" = 30;
    "E15" = 175;
    "F15" = 99;
    "H15" = 140;
    "E17" = 136;
    "E18" = 131;
    "F18" = 61;
    "G18" = 37;
    "H18" = 98;
    "E24" = 28;
    "F28" = 19;
    "H28" = 21;
    "E41" = 50;
    "E44" = 32;
    "E48" = 40;
    "F48" = 26;
    "H48" = 32;
    "E60" = 21;
    "E61" = 36;
    "F61" = 16;
    "H61" = 26;
    "E62" = 54;
    "E63" = 47;
    "F65" = 15;
    "H65" = 28;
    "F66" = 27;
    "H66" = 35;
    "E71" = 44;
    "E75" = 20;
    "E76" = 57;
    "E79" = 42;
    "F79" = 20;
    "H79" = 31;
    "E80" = 34;
    "F80" = 18;
    "H80" = 30;
    "E82" = 19;
    "F82" = 8;
    "H82" = 14;
    "E86" = 5;
    "F86" = 3;
    "H86" = 5;
    "F88" = 19;
    "H88" = 27;
    "E89" = 48;
    "G89" = 6;
    "H89" = 29;
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
